# Generate Report for Handoff
# Updates status text + handoff timestamps to reflect that the report has
# moved from "In Translation" to "Ready for handoff", and widens the
# "Status" columns to fit the new (longer) text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
# E2 / F2: per-locale status; G2: latest handoff-xliff-generate datetime
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-31 14:49:30"

# --- zh-cn sheet -------------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-31 14:49:24"

# --- de-de sheet -------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-31 14:49:30"

# --- Column widths ------------------------------------------------------
# The Status columns grew to fit "Ready for handoff" (longer than
# "In Translation").
$overview.Range("E1").ColumnWidth = 17.2159881591797
$overview.Range("F1").ColumnWidth = 17.2159881591797
$zhcn.Range("C1").ColumnWidth = 17.2159881591797
$dede.Range("C1").ColumnWidth = 17.2159881591797
